$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Multi-Utilities(18)'
$ws.Cells.Item(2, 2).Value = 0.4701696426296317

$ws.Cells.Item(3, 1).Value = 'Energy Equipment & Services(32)'
$ws.Cells.Item(3, 2).Value = 0.406557305300907

$ws.Cells.Item(4, 1).Value = 'Road & Rail(22)'
$ws.Cells.Item(4, 2).Value = 0.3673881964468256

$ws.Cells.Item(5, 1).Value = 'Banks(246)'
$ws.Cells.Item(5, 2).Value = 0.3368748846693499

$ws.Cells.Item(6, 1).Value = 'Electric Utilities(28)'
$ws.Cells.Item(6, 2).Value = 0.3205586396304085

$ws.Cells.Item(7, 1).Value = 'Marine(15)'
$ws.Cells.Item(7, 2).Value = 0.3070128247586821

$ws.Cells.Item(8, 1).Value = 'Building Products(23)'
$ws.Cells.Item(8, 2).Value = 0.302143338175924

$ws.Cells.Item(9, 1).Value = 'Auto Components(21)'
$ws.Cells.Item(9, 2).Value = 0.2938848054147093

$ws.Cells.Item(10, 1).Value = 'Machinery(85)'
$ws.Cells.Item(10, 2).Value = 0.2664101813557964

$ws.Cells.Item(11, 1).Value = 'Construction & Engineering(20)'
$ws.Cells.Item(11, 2).Value = 0.2616263087945421

$ws.Cells.Item(12, 1).Value = 'Trading Companies & Distributors(25)'
$ws.Cells.Item(12, 2).Value = 0.2552831537526278

$ws.Cells.Item(13, 1).Value = 'Specialty Retail(58)'
$ws.Cells.Item(13, 2).Value = 0.2498051986964366

$ws.Cells.Item(14, 1).Value = 'Thrifts & Mortgage Finance(47)'
$ws.Cells.Item(14, 2).Value = 0.2357986428000619

$ws.Cells.Item(15, 1).Value = 'Textiles, Apparel & Luxury Goods(29)'
$ws.Cells.Item(15, 2).Value = 0.2332191071178296

$ws.Cells.Item(16, 1).Value = 'Capital Markets(75)'
$ws.Cells.Item(16, 2).Value = 0.1909815052457282

$ws.Cells.Item(17, 1).Value = 'Hotels, Restaurants & Leisure(50)'
$ws.Cells.Item(17, 2).Value = 0.1870392878746558

$ws.Cells.Item(18, 1).Value = 'Insurance(75)'
$ws.Cells.Item(18, 2).Value = 0.1840892825553289

$ws.Cells.Item(19, 1).Value = 'Oil, Gas & Consumable Fuels(122)'
$ws.Cells.Item(19, 2).Value = 0.1837975971606767

$ws.Cells.Item(20, 1).Value = 'Semiconductors & Semiconductor Equipment(68)'
$ws.Cells.Item(20, 2).Value = 0.1803422401853819

$ws.Cells.Item(21, 1).Value = 'Professional Services(35)'
$ws.Cells.Item(21, 2).Value = 0.1672019193938927

$ws.Cells.Item(22, 1).Value = 'IT Services(52)'
$ws.Cells.Item(22, 2).Value = 0.1491692899065586

$ws.Cells.Item(23, 1).Value = 'Chemicals(51)'
$ws.Cells.Item(23, 2).Value = 0.1458712031784943

$ws.Cells.Item(24, 1).Value = 'Metals & Mining(89)'
$ws.Cells.Item(24, 2).Value = 0.1440840363727988

$ws.Cells.Item(25, 1).Value = 'Household Durables(39)'
$ws.Cells.Item(25, 2).Value = 0.1340439084648633

$ws.Cells.Item(26, 1).Value = 'Aerospace & Defense(37)'
$ws.Cells.Item(26, 2).Value = 0.1328854393022718

$ws.Cells.Item(27, 1).Value = 'Health Care Equipment & Supplies(83)'
$ws.Cells.Item(27, 2).Value = 0.1205744696534682

$ws.Cells.Item(28, 1).Value = 'Health Care Providers & Services(46)'
$ws.Cells.Item(28, 2).Value = 0.1164692854126994

$ws.Cells.Item(29, 1).Value = 'Commercial Services & Supplies(52)'
$ws.Cells.Item(29, 2).Value = 0.1076821198739539

$ws.Cells.Item(30, 1).Value = 'Software(66)'
$ws.Cells.Item(30, 2).Value = 0.08182842459818392

$ws.Cells.Item(31, 1).Value = 'Biotechnology(126)'
$ws.Cells.Item(31, 2).Value = 0.07576123345670216

